$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark from around the
#    "Additional passes may be acquired..." paragraph to wrap the
#    "Any vehicle parked for longer than 20 minutes..." paragraph.
# ------------------------------------------------------------------
$anyVehicleRange = $d.Content
$anyVehicleRange.Find.Execute("Any vehicle parked for longer than 20 minutes will be cited or towed away at owner's expense.") | Out-Null
$anyVehicleParagraph = $anyVehicleRange.Paragraphs(1)

# Remove the pre-existing _GoBack bookmark (currently sitting at the
# end of the "Additional passes..." paragraph).
$goBackCount = 0
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $goBackCount = $goBackCount + 1
    }
}
if ($goBackCount -gt 0) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create it spanning the "Any vehicle parked..." paragraph.
$d.Bookmarks.Add("_GoBack", $anyVehicleParagraph.Range) | Out-Null

# ------------------------------------------------------------------
# 2. Add four new empty, bold-formatted paragraphs right after the
#    "Additional passes may be acquired..." paragraph.
# ------------------------------------------------------------------
$additionalPassesRange = $d.Content
$additionalPassesRange.Find.Execute("Additional passes may be acquired if needed at the check-in desk as space permits.") | Out-Null
$additionalPassesParagraph = $additionalPassesRange.Paragraphs(1)

$insertionPoint = $additionalPassesParagraph.Range
for ($i = 0; $i -lt 4; $i++) {
    $insertionPoint.InsertParagraphAfter() | Out-Null
    $insertionPoint = $insertionPoint.Next(4) # wdParagraph = 4
}

Write-Host "done"
